$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping from target row (2..17) to the original source row (2..19) that
# contains the data which must end up there. Rows 3 ("스튜디오삼익") and
# 10 ("신영스팩10호") of the original sheet are dropped (their entity, 신영,
# is no longer present in the dataset). The remaining rows are reordered as
# shown below (some blocks are cyclically rotated in the source data).
$srcRows = @(2,4,5,7,8,6,9,11,12,13,14,17,15,16,18,19)

# Stage the source rows (2:19) into a scratch area first (copy, not move),
# so that writes to the destination area don't clobber rows we still need
# to read from later in the loop.
$ws.Range("A2:L19").Copy($ws.Range("A100:L117"))

for ($i = 0; $i -lt $srcRows.Length; $i++) {
    $destRow = $i + 2
    $scratchRow = $srcRows[$i] + 98
    $ws.Range("A" + $scratchRow + ":L" + $scratchRow).Copy($ws.Range("A" + $destRow + ":L" + $destRow))
}

# Remove the scratch area used for staging.
$ws.Range("A100:L117").Clear()

# Remove the now-superfluous trailing rows (18 and 19) left over from the
# original 19-row table; the table now only spans through row 17.
$ws.Rows("18:19").Delete()
